$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 4, shifting everything else down by one.
$ws.Rows("4:4").Insert()

# New row 4 becomes a link to the Project Euler problem page.
$ws.Hyperlinks.Add($ws.Range("B4"), "https://projecteuler.net/problem=5")

# Move the view back to the top and select F3 (matches the saved view state).
[void]$ws.Range("F3").Select()

# Italicise the two explanatory lines that used to be rows 4-5 (now rows 5-6).
$ws.Range("B5").Font.Italic = $true
$ws.Range("B6").Font.Italic = $true

# Within the rich-text question cell (now B6), italicise the two trailing runs
# ("evenly divisible" and the remainder of the sentence) to match the target
# formatting while leaving the leading run untouched.
$cell = $ws.Range("B6")
$cell.Characters(46, 16).Font.Italic = $true
$cell.Characters(62, 36).Font.Italic = $true
